# This script reorders (permutes) the per-row values of columns
# D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado) and P (Precio $/Kg) across rows 2-31.
# Each destination row receives the tuple of values that used to belong
# to a (different) source row, while all other columns (A,B,C,E,F,G,H,I,
# N,O,Q,R) remain untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (values are read from the
# source row's original contents and written into the destination row).
$map = @{
    2  = 21
    3  = 16
    4  = 6
    5  = 12
    6  = 14
    7  = 13
    8  = 7
    9  = 29
    10 = 31
    11 = 15
    12 = 26
    13 = 24
    14 = 22
    15 = 28
    16 = 10
    17 = 5
    18 = 9
    19 = 19
    20 = 8
    21 = 3
    22 = 20
    23 = 27
    24 = 11
    25 = 2
    26 = 23
    27 = 17
    28 = 4
    29 = 25
    30 = 30
    31 = 18
}

$cols = @("D", "J", "K", "L", "M", "P")

# Snapshot all original values first so reads are never affected by
# writes that happen earlier in the loop. Value2 is used for reading
# because it returns the raw number reliably in this environment.
$original = @{}
foreach ($row in $map.Keys) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$row").Value2
    }
    $original[$row] = $rowVals
}

foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    $srcVals = $original[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $srcVals[$col]
    }
}
